$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '24.622.82'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.20%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.689.09'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.42'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.60%  '

$ws.Range('E6').Value = '  +0.09%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3898'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.97%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4030'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.47%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.496'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.81%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.004'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.12%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.77'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.31%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08740'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.04%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.598'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.02%  '

$ws.Range('E14').Value = '  +5.68%  '

$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001350'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.56%  '

$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.958'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.70%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.685.18'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.18%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '98.33'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.23%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07095'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.22%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.76'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.37%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.279'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.24%  '

$ws.Range('E22').Value = '  -0.35%  '

$ws.Range('E23').Value = '  -0.48%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '24.627.04'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.10%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.003'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -8.66%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.352'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.12%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.70'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.04%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '161.49'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.73%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.596'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +13.02%  '

$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '136.23'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.73%  '

$ws.Range('B31').Value = 'HuobiToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.212'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.66%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.869.83'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.25%  '

$ws.Range('E33').Value = '  +2.76%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.389'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.05%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.039'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.84%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.987'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.39%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02907'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.14%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2716'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.36%  '

$ws.Range('E39').Value = '  -4.46%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '14.19'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.63%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.09111'
$ws.Range('D41').Style = 'Normal'

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7798'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.32%  '

$ws.Range('E43').Value = '  -0.61%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.59'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.99%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.7188'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.83%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.578'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.34%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.193'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.49%  '

$ws.Range('E48').Value = '  +0.03%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.334'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.13%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '137.65'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.31%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '90.80'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.31%  '
